$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7220110893249512
$ws.Range("E2").Value = 1823.946026437212
$ws.Range("G2").Value = 0.09755426506924006
$ws.Range("H2").Value = 0.07550795654597584
$ws.Range("I2").Value = 0.06697381020525341
$ws.Range("J2").Value = 0.05987593863057437
$ws.Range("K2").Value = 0.05574161441637548
$ws.Range("L2").Value = 0.04823549533296802
$ws.Range("M2").Value = 0.04595180715092104
$ws.Range("N2").Value = 0.04441975026484835
$ws.Range("O2").Value = 0.04276830668710583
$ws.Range("P2").Value = 0.04080214802851328
$ws.Range("Q2").Value = 0.04000859774249649
$ws.Range("R2").Value = 0.03877391321004008
$ws.Range("S2").Value = 0.03858781520965183
$ws.Range("T2").Value = 0.0374887958356401
$ws.Range("U2").Value = 0.0369286922833125
$ws.Range("V2").Value = 0.03672369529029165
$ws.Range("W2").Value = 0.03616662572769774
$ws.Range("X2").Value = 0.03599851515851746
$ws.Range("Y2").Value = 0.03555450343932186
$ws.Range("C3").Value = 0.8376443386077881
$ws.Range("E3").Value = 1920.719378025262
$ws.Range("G3").Value = 0.09607233796459398
$ws.Range("H3").Value = 0.08433468883053087
$ws.Range("I3").Value = 0.07376400577307486
$ws.Range("J3").Value = 0.06089228016560137
$ws.Range("K3").Value = 0.05560544506762369
$ws.Range("L3").Value = 0.05073648509466843
$ws.Range("M3").Value = 0.04898786712668331
$ws.Range("N3").Value = 0.04535042575653891
$ws.Range("O3").Value = 0.04351389479441652
$ws.Range("P3").Value = 0.04163223637581933
$ws.Range("Q3").Value = 0.04141649103004245
$ws.Range("R3").Value = 0.03982333529340368
$ws.Range("S3").Value = 0.03982333529340368
$ws.Range("T3").Value = 0.03898243912339753
$ws.Range("U3").Value = 0.0383735856373338
$ws.Range("V3").Value = 0.03809145371784466
$ws.Range("W3").Value = 0.03781107411420781
$ws.Range("X3").Value = 0.03749893669352677
$ws.Range("Y3").Value = 0.03744092354825071
$ws.Range("C4").Value = 0.7041990756988525
$ws.Range("E4").Value = 1886.015702659501
$ws.Range("G4").Value = 0.1037707458911285
$ws.Range("H4").Value = 0.07992326389406258
$ws.Range("I4").Value = 0.06791069735184614
$ws.Range("J4").Value = 0.06163287338387032
$ws.Range("K4").Value = 0.05591880462104015
$ws.Range("L4").Value = 0.05220236302705725
$ws.Range("M4").Value = 0.04951235336616199
$ws.Range("N4").Value = 0.04667236603712394
$ws.Range("O4").Value = 0.04471364762593191
$ws.Range("P4").Value = 0.04278368768147681
$ws.Range("Q4").Value = 0.04218357575354013
$ws.Range("R4").Value = 0.03932815765306952
$ws.Range("S4").Value = 0.03883567549070741
$ws.Range("T4").Value = 0.03818919790285589
$ws.Range("U4").Value = 0.03762709490733979
$ws.Range("V4").Value = 0.03726176753103055
$ws.Range("W4").Value = 0.03702184489040632
$ws.Range("X4").Value = 0.0368001986716062
$ws.Range("Y4").Value = 0.03676443864833333
$ws.Range("C5").Value = 0.7344028949737549
$ws.Range("E5").Value = 1931.105843338913
$ws.Range("G5").Value = 0.09972497870805652
$ws.Range("H5").Value = 0.07513082189433373
$ws.Range("I5").Value = 0.06854983708743935
$ws.Range("J5").Value = 0.06243985869062176
$ws.Range("K5").Value = 0.05725564188087752
$ws.Range("L5").Value = 0.05168611458303569
$ws.Range("M5").Value = 0.04869201823428484
$ws.Range("N5").Value = 0.0462414341308875
$ws.Range("O5").Value = 0.04383203667808887
$ws.Range("P5").Value = 0.04208470042088369
$ws.Range("Q5").Value = 0.04095163990377351
$ws.Range("R5").Value = 0.03961527453832409
$ws.Range("S5").Value = 0.03831237853053068
$ws.Range("T5").Value = 0.03831237853053068
$ws.Range("U5").Value = 0.03831237853053068
$ws.Range("V5").Value = 0.03826520469065586
$ws.Range("W5").Value = 0.03808250152023961
$ws.Range("X5").Value = 0.03781917945394159
$ws.Range("Y5").Value = 0.03764338875904314
$ws.Range("C6").Value = 0.7187607288360596
$ws.Range("E6").Value = 1884.856715977297
$ws.Range("G6").Value = 0.1014846073387392
$ws.Range("H6").Value = 0.08414089153894179
$ws.Range("I6").Value = 0.07123269455373123
$ws.Range("J6").Value = 0.06055507959823386
$ws.Range("K6").Value = 0.05499682706365457
$ws.Range("L6").Value = 0.05108291196578486
$ws.Range("M6").Value = 0.04760584005438271
$ws.Range("N6").Value = 0.04477477974950918
$ws.Range("O6").Value = 0.04277751515394787
$ws.Range("P6").Value = 0.04118777923749412
$ws.Range("Q6").Value = 0.04059687347772081
$ws.Range("R6").Value = 0.03993539287764056
$ws.Range("S6").Value = 0.03916631684820188
$ws.Range("T6").Value = 0.03812719086749505
$ws.Range("U6").Value = 0.03746007862497041
$ws.Range("V6").Value = 0.03739243728132347
$ws.Range("W6").Value = 0.03722061967454487
$ws.Range("X6").Value = 0.03702335885966018
$ws.Range("Y6").Value = 0.03674184631534692
$ws.Range("C7").Value = 0.7187397480010986
$ws.Range("E7").Value = 1854.940983446664
$ws.Range("G7").Value = 0.1089045159863328
$ws.Range("H7").Value = 0.08421659745009211
$ws.Range("I7").Value = 0.0700224624433204
$ws.Range("J7").Value = 0.06279322913652802
$ws.Range("K7").Value = 0.05779449765231579
$ws.Range("L7").Value = 0.05137961801488286
$ws.Range("M7").Value = 0.0494111291254624
$ws.Range("N7").Value = 0.04601578470339477
$ws.Range("O7").Value = 0.04362483338075326
$ws.Range("P7").Value = 0.0419397999239084
$ws.Range("Q7").Value = 0.04045012649036446
$ws.Range("R7").Value = 0.03906696264959705
$ws.Range("S7").Value = 0.03882425573428665
$ws.Range("T7").Value = 0.03795734671708675
$ws.Range("U7").Value = 0.03761693130795402
$ws.Range("V7").Value = 0.03706481857656186
$ws.Range("W7").Value = 0.03669644835510839
$ws.Range("X7").Value = 0.03640554096039761
$ws.Range("Y7").Value = 0.03615869363443789
$ws.Range("C8").Value = 0.7031242847442627
$ws.Range("E8").Value = 1905.348472241387
$ws.Range("G8").Value = 0.1007036288752944
$ws.Range("H8").Value = 0.07946318827017995
$ws.Range("I8").Value = 0.07009450254034896
$ws.Range("J8").Value = 0.05943356081180183
$ws.Range("K8").Value = 0.05669277646733819
$ws.Range("L8").Value = 0.05279198122545565
$ws.Range("M8").Value = 0.04971447293176867
$ws.Range("N8").Value = 0.04728085800719274
$ws.Range("O8").Value = 0.04352446006136094
$ws.Range("P8").Value = 0.04282995581924223
$ws.Range("Q8").Value = 0.04142241657266323
$ws.Range("R8").Value = 0.04044523446675234
$ws.Range("S8").Value = 0.03889661795742109
$ws.Range("T8").Value = 0.03820346730827238
$ws.Range("U8").Value = 0.03800596266885135
$ws.Range("V8").Value = 0.03758576288017162
$ws.Range("W8").Value = 0.03744597526049239
$ws.Range("X8").Value = 0.03733795134652403
$ws.Range("Y8").Value = 0.03714129575519272
$ws.Range("C9").Value = 0.703136682510376
$ws.Range("E9").Value = 1935.458128910157
$ws.Range("G9").Value = 0.09534914520707652
$ws.Range("H9").Value = 0.07224515919959679
$ws.Range("I9").Value = 0.06623905826503551
$ws.Range("J9").Value = 0.05723152497131586
$ws.Range("K9").Value = 0.05499973095315536
$ws.Range("L9").Value = 0.05165286178759043
$ws.Range("M9").Value = 0.0484177578299974
$ws.Range("N9").Value = 0.04622649038476335
$ws.Range("O9").Value = 0.04456887353344412
$ws.Range("P9").Value = 0.0433201527890749
$ws.Range("Q9").Value = 0.04270557937040182
$ws.Range("R9").Value = 0.04135898611660306
$ws.Range("S9").Value = 0.04041427346188384
$ws.Range("T9").Value = 0.03956561781697337
$ws.Range("U9").Value = 0.03900476392637871
$ws.Range("V9").Value = 0.03855643669713028
$ws.Range("W9").Value = 0.03807897066610921
$ws.Range("X9").Value = 0.03787559846842762
$ws.Range("Y9").Value = 0.03772822863372625
$ws.Range("C10").Value = 0.7187502384185791
$ws.Range("E10").Value = 1838.21125439734
$ws.Range("G10").Value = 0.1058780337508392
$ws.Range("H10").Value = 0.08210602197036311
$ws.Range("I10").Value = 0.06794050498999592
$ws.Range("J10").Value = 0.05918352216677716
$ws.Range("K10").Value = 0.05538444670318474
$ws.Range("L10").Value = 0.04927723239472025
$ws.Range("M10").Value = 0.04792240452736134
$ws.Range("N10").Value = 0.04597571676709944
$ws.Range("O10").Value = 0.04262993203026734
$ws.Range("P10").Value = 0.0405487358257329
$ws.Range("Q10").Value = 0.0405487358257329
$ws.Range("R10").Value = 0.03949815950609844
$ws.Range("S10").Value = 0.03823199609927602
$ws.Range("T10").Value = 0.03751166023074618
$ws.Range("U10").Value = 0.03669820422602187
$ws.Range("V10").Value = 0.03653204022592882
$ws.Range("W10").Value = 0.03630197082676493
$ws.Range("X10").Value = 0.03616521001987217
$ws.Range("Y10").Value = 0.03583257805842767
$ws.Range("C11").Value = 0.7031130790710449
$ws.Range("E11").Value = 1905.831279267479
$ws.Range("G11").Value = 0.1003861996418925
$ws.Range("H11").Value = 0.07384283783350278
$ws.Range("I11").Value = 0.06966025215566361
$ws.Range("J11").Value = 0.06407893081725483
$ws.Range("K11").Value = 0.05962623123572769
$ws.Range("L11").Value = 0.05303031933438934
$ws.Range("M11").Value = 0.05137566469027309
$ws.Range("N11").Value = 0.04725832973786386
$ws.Range("O11").Value = 0.04516137433891549
$ws.Range("P11").Value = 0.04213712632814529
$ws.Range("Q11").Value = 0.04178215554414532
$ws.Range("R11").Value = 0.04035641349879685
$ws.Range("S11").Value = 0.03949328053575445
$ws.Range("T11").Value = 0.03876250710612907
$ws.Range("U11").Value = 0.03807792600262717
$ws.Range("V11").Value = 0.03803680678635471
$ws.Range("W11").Value = 0.0375923985127471
$ws.Range("X11").Value = 0.03729040252534555
$ws.Range("Y11").Value = 0.03715070719819646
